# Update stats after logging the 2021 divisional round game,
# simulated season continuing from the conference round.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 4: A.Dillon
$rushing.Range("C4").Value = 118
$rushing.Range("D4").Value = 59
$rushing.Range("E4").Value = 11
$rushing.Range("F4").Value = 35

# Row 5: K.Hill
$rushing.Range("C5").Value = 90
$rushing.Range("D5").Value = 65
$rushing.Range("F5").Value = 38

# Row 10: D.Dafney
$rushing.Range("C10").Value = 4

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: A.Dillon
$receiving.Range("C2").Value = 73
$receiving.Range("D2").Value = 61
$receiving.Range("E2").Value = 7
$receiving.Range("F2").Value = 5
$receiving.Range("G2").Value = 17
$receiving.Range("H2").Value = 14

# Row 5: D.Adams
$receiving.Range("C5").Value = 142
$receiving.Range("D5").Value = 118
$receiving.Range("E5").Value = 42
$receiving.Range("F5").Value = 31
$receiving.Range("G5").Value = 30
$receiving.Range("H5").Value = 25

# Row 7: A.Lazard
$receiving.Range("C7").Value = 51
$receiving.Range("D7").Value = 40
$receiving.Range("G7").Value = 15
$receiving.Range("H7").Value = 8

# Row 8: M.Valdes-Scantling
$receiving.Range("C8").Value = 34

# Row 13: M.Lewis
$receiving.Range("C13").Value = 27
$receiving.Range("D13").Value = 24

# Row 14: J.Deguara
$receiving.Range("E14").Value = 4

# Row 15: D.Dafney
$receiving.Range("C15").Value = 3
$receiving.Range("E15").Value = 1
